$d = $word.ActiveDocument

$d.Content.Find.Execute("[[PERSON_11]] – „s [[PERSON_12]]“, „o [[PERSON_13]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_11]] – „s [[PERSON_12]]“, „o [[PERSON_11]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_14]] – „ke [[PERSON_14]]“, „o [[PERSON_14]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_13]] – „ke [[PERSON_13]]“, „o [[PERSON_13]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_15]] – „o [[PERSON_16]]“, „s [[PERSON_15]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_14]] – „o [[PERSON_15]]“, „s [[PERSON_14]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_17]] – „u [[PERSON_17]]“, „s [[PERSON_17]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_16]] – „u [[PERSON_16]]“, „s [[PERSON_16]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_18]] – „ke [[PERSON_18]]“, „o [[PERSON_18]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_17]] – „ke [[PERSON_17]]“, „o [[PERSON_17]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_19]] – „s [[PERSON_19]]“, „o [[PERSON_19]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_18]] – „s [[PERSON_18]]“, „o [[PERSON_18]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_20]] – „s [[PERSON_20]]“, „o [[PERSON_20]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_19]] – „s [[PERSON_19]]“, „o [[PERSON_19]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_21]] – „k [[PERSON_21]]“, „od [[PERSON_21]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_20]] – „k [[PERSON_20]]“, „od [[PERSON_20]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_22]] – „o [[PERSON_22]]“, „s [[PERSON_22]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_21]] – „o [[PERSON_21]]“, „s [[PERSON_21]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_23]] – „o [[PERSON_24]]“, „se [[PERSON_25]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_22]] – „o [[PERSON_23]]“, „se [[PERSON_24]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_26]] – „s [[PERSON_26]]“, „u [[PERSON_26]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_25]] – „s [[PERSON_25]]“, „u [[PERSON_25]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_27]] – „o [[PERSON_28]]“, „s [[PERSON_29]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_26]] – „o [[PERSON_26]]“, „s [[PERSON_27]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_30]] – „k [[PERSON_30]]“, „o [[PERSON_30]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_28]] – „k [[PERSON_28]]“, „o [[PERSON_28]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_31]] – „se [[PERSON_31]]“, „o Soně Mikulkové“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_29]] – „se [[PERSON_29]]“, „o Soně Mikulkové“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_32]] – „o [[PERSON_32]]“, „s [[PERSON_32]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_30]] – „o [[PERSON_30]]“, „s [[PERSON_30]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_33]] – „s [[PERSON_33]]“, „o [[PERSON_33]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_31]] – „s [[PERSON_31]]“, „o [[PERSON_31]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_34]] – „k [[PERSON_35]]“, „s [[PERSON_34]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_32]] – „k [[PERSON_33]]“, „s [[PERSON_32]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_36]] – „s [[PERSON_36]]“, „o [[PERSON_37]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_34]] – „s [[PERSON_34]]“, „o [[PERSON_35]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_38]] – „od [[PERSON_38]]“, „s [[PERSON_39]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_36]] – „od [[PERSON_36]]“, „s [[PERSON_37]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_40]] – „k [[PERSON_41]]“, „o [[PERSON_41]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_38]] – „k [[PERSON_39]]“, „o [[PERSON_39]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_42]] – „o [[PERSON_43]]“, „s [[PERSON_42]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_40]] – „o [[PERSON_41]]“, „s [[PERSON_40]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_44]] – „s [[PERSON_44]]“, „o [[PERSON_45]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_42]] – „s [[PERSON_42]]“, „o [[PERSON_43]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_47]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_44]] – „s [[PERSON_44]]“, „o [[PERSON_45]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_48]] – „k [[PERSON_48]]“, „s [[PERSON_48]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_46]] – „k [[PERSON_46]]“, „s [[PERSON_46]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_49]] – „pro [[PERSON_50]]“, „o [[PERSON_51]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_47]] – „pro [[PERSON_48]]“, „o [[PERSON_49]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_52]] – „k [[PERSON_52]]“, „o [[PERSON_52]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_50]] – „k [[PERSON_50]]“, „o [[PERSON_50]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_53]] – „o [[PERSON_54]]“, „s [[PERSON_53]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_51]] – „o [[PERSON_52]]“, „s [[PERSON_51]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_55]] – „s [[PERSON_55]]“, „o [[PERSON_56]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_53]] – „s [[PERSON_53]]“, „o [[PERSON_54]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_57]] – „s [[PERSON_57]]“, „o [[PERSON_57]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_55]] – „s [[PERSON_55]]“, „o [[PERSON_55]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_58]] – „u [[PERSON_58]]“, „o [[PERSON_59]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_56]] – „u [[PERSON_56]]“, „o [[PERSON_57]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_60]] – „se [[PERSON_60]]“, „o [[PERSON_60]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_58]] – „se [[PERSON_58]]“, „o [[PERSON_58]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_61]] – „o [[PERSON_62]]“, „s [[PERSON_63]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_59]] – „o [[PERSON_60]]“, „s [[PERSON_61]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_64]] – „k [[PERSON_65]]“, „o [[PERSON_65]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_62]] – „k [[PERSON_63]]“, „o [[PERSON_63]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_66]] – „o [[PERSON_67]]“, „s [[PERSON_66]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_64]] – „o [[PERSON_65]]“, „s [[PERSON_64]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_68]] – „s [[PERSON_68]]“, „o [[PERSON_68]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_66]] – „s [[PERSON_66]]“, „o [[PERSON_66]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_69]] – „s [[PERSON_69]]“, „o [[PERSON_70]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_67]] – „s [[PERSON_67]]“, „o [[PERSON_68]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_71]] – „o [[PERSON_71]]“, „s [[PERSON_71]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_69]] – „o [[PERSON_69]]“, „s [[PERSON_69]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_72]] – „s [[PERSON_73]]“, „o [[PERSON_74]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_70]] – „s [[PERSON_71]]“, „o [[PERSON_72]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_75]] – „o [[PERSON_76]]“, „s [[PERSON_75]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_73]] – „o [[PERSON_74]]“, „s [[PERSON_73]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_77]] – „s [[PERSON_78]]“, „o [[PERSON_79]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_75]] – „s [[PERSON_76]]“, „o [[PERSON_77]]“", 2) | Out-Null
